$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(57).Insert()

$ws.Cells.Item(57, 1).Value = 11
$ws.Cells.Item(57, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(57, 3).Value = "Bíobío"
$ws.Cells.Item(57, 4).Value = (Get-Date -Year 2021 -Month 9 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(57, 5).Value = 8
$ws.Cells.Item(57, 6).Value = "Fruta"
$ws.Cells.Item(57, 7).Value = 100101
$ws.Cells.Item(57, 8).Value = "Berries"
$ws.Cells.Item(57, 9).Value = 100112025
$ws.Cells.Item(57, 10).Value = "Frutilla"
$ws.Cells.Item(57, 11).Value = "Sin especificar"
$ws.Cells.Item(57, 12).Value = "Primera"
$ws.Cells.Item(57, 13).Value = 50
$ws.Cells.Item(57, 14).Value = 17000
$ws.Cells.Item(57, 15).Value = 17000
$ws.Cells.Item(57, 16).Value = 17000
$ws.Cells.Item(57, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(57, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(57, 19).Value = 2429
$ws.Cells.Item(57, 20).Value = 7
